$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("D2").Value = "67.757.49"
$ws.Range("E2").Value = "  +0.28%  "
$ws.Range("D3").Value = "3.815.87"
$ws.Range("E3").Value = "  +1.34%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "601.78"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.25%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "166.14"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.19%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.22%  "
$ws.Range("E8").Value = "  +0.06%  "
$ws.Range("E9").Value = "  +0.63%  "
$ws.Range("E10").Value = "  -0.91%  "
$ws.Range("E11").Value = "  +0.54%  "
$ws.Range("E12").Value = "  -0.75%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "35.78"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.84%  "
$ws.Range("D14").Value = "4.456.69"
$ws.Range("E14").Value = "  +1.36%  "
$ws.Range("D15").Value = "3.805.22"
$ws.Range("E15").Value = "  +1.59%  "
$ws.Range("D16").Value = "67.776.56"
$ws.Range("E16").Value = "  +0.38%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "18.43"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.05%  "
$ws.Range("E18").Value = "  +1.02%  "
$ws.Range("E19").Value = "  +0.35%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "463.33"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.77%  "
$ws.Range("E21").Value = "  -0.95%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.699"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.51%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.0000149"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -3.86%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "83.41"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.15%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "12.09"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.91%  "
$ws.Range("E26").Value = "  -1.64%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.05"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.88%  "
$ws.Range("E28").Value = "  -0.08%  "
$ws.Range("D29").Value = "3.965.83"
$ws.Range("E29").Value = "  +1.40%  "
$ws.Range("E30").Value = "  +0.33%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.37"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.61%  "
$ws.Range("E32").Value = "  +1.87%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "29.54"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.35%  "
$ws.Range("E34").Value = "  +0.16%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "9.08"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.82%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.1000"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.30%  "
$ws.Range("B37").Value = "dogwifhat"
$ws.Range("C37").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.28"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.80%  "
$ws.Range("B38").Value = "Kaspa"
$ws.Range("C38").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.138"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.12%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.999"
$ws.Range("D39").Style = "Normal"
$ws.Range("E40").Value = "  +0.90%  "
$ws.Range("E41").Value = "  +0.02%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "48.07"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.23%  "
$ws.Range("E44").Value = "  +0.48%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "28.45"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +9.60%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "43.13"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -4.90%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.41"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +13.31%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "148.91"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.19%  "
$ws.Range("E49").Value = "  +0.19%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.84"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.15%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "386.19"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.91%  "
